# Scheduled-runner price refresh: rewrite computed profit columns (H:N) on the
# affected leve rows across all eight job sheets. Values are literal numbers
# (no formulas in this workbook), matching how the runner writes results back.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @(80, "H", 631),
    @(80, "I", 555.1667),
    @(80, "J", 858.5),
    @(80, "K", 1665.5001),
    @(80, "L", 2575.5),
    @(80, "M", -667.5001),
    @(80, "N", -4571.5),
    @(83, "H", 631),
    @(83, "I", 555.1667),
    @(83, "J", 858.5),
    @(83, "K", 4996.5003),
    @(83, "L", 7726.5),
    @(83, "M", -4.500299999999697),
    @(83, "N", -17710.5),
    @(132, "H", 1259.25),
    @(132, "I", 1111.4),
    @(132, "K", 3334.2),
    @(132, "M", -804.2000000000003),
    @(135, "H", 4646.4),
    @(135, "I", 5558),
    @(135, "J", 1000),
    @(135, "K", 50022),
    @(135, "L", 9000),
    @(135, "M", -47487),
    @(135, "N", -14070),
    @(137, "H", 3476.1538),
    @(137, "I", 741),
    @(137, "K", 2223),
    @(137, "M", 327),
    @(138, "H", 3948.2942),
    @(138, "I", 3126.2),
    @(138, "J", 4090.0344),
    @(138, "K", 9378.599999999999),
    @(138, "L", 12270.1032),
    @(138, "M", -4238.599999999999),
    @(138, "N", -22550.1032),
    @(141, "H", 7692.5557),
    @(141, "I", 7798.8),
    @(141, "J", 7559.75),
    @(141, "K", 23396.4),
    @(141, "L", 22679.25),
    @(141, "M", -18216.4),
    @(141, "N", -33039.25)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @(61, "H", 2905.2),
    @(61, "I", 2783.7778),
    @(61, "K", 2783.7778),
    @(61, "M", -2571.7778),
    @(74, "H", 11576.75),
    @(74, "I", 7019.1665),
    @(74, "J", 25249.5),
    @(74, "K", 7019.1665),
    @(74, "L", 25249.5),
    @(74, "M", -6145.1665),
    @(74, "N", -26997.5),
    @(77, "H", 11576.75),
    @(77, "I", 7019.1665),
    @(77, "J", 25249.5),
    @(77, "K", 35095.8325),
    @(77, "L", 126247.5),
    @(77, "M", -30727.8325),
    @(77, "N", -134983.5),
    @(136, "H", 2905.2),
    @(136, "I", 2783.7778),
    @(136, "K", 8351.3334),
    @(136, "M", -5801.3334)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @(17, "H", 16504),
    @(17, "J", 16504),
    @(17, "L", 16504),
    @(17, "N", -16852),
    @(25, "H", 15000),
    @(25, "I", 15000),
    @(25, "K", 15000),
    @(25, "M", -14826),
    @(31, "H", 6262.36),
    @(31, "I", 2572.0715),
    @(31, "K", 2572.0715),
    @(31, "M", -2277.0715),
    @(34, "H", 6262.36),
    @(34, "I", 2572.0715),
    @(34, "K", 2572.0715),
    @(34, "M", -2370.0715),
    @(58, "H", 2191.818),
    @(58, "I", 2211),
    @(58, "J", 2000),
    @(58, "K", 2211),
    @(58, "L", 2000),
    @(58, "M", -2008),
    @(58, "N", -2406),
    @(132, "H", 3524.318),
    @(132, "I", 2855.2354),
    @(132, "K", 8565.706200000001),
    @(132, "M", -6035.706200000001),
    @(134, "H", 6599.4),
    @(134, "I", 6599.4),
    @(134, "J", 0),
    @(134, "K", 19798.2),
    @(134, "L", 0),
    @(134, "M", -17263.2),
    @(134, "N", $null),
    @(136, "H", 2191.818),
    @(136, "I", 2211),
    @(136, "J", 2000),
    @(136, "K", 6633),
    @(136, "L", 6000),
    @(136, "M", -4083),
    @(136, "N", -11100)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @(5, "H", 800.8),
    @(5, "I", 701.5),
    @(5, "J", 867),
    @(5, "K", 2104.5),
    @(5, "L", 2601),
    @(5, "M", -1992.5),
    @(5, "N", -2825),
    @(9, "H", 44400.285),
    @(9, "I", 1300),
    @(9, "J", 51583.668),
    @(9, "K", 3900),
    @(9, "L", 154751.004),
    @(9, "M", -3676),
    @(9, "N", -155199.004),
    @(11, "H", 712.25),
    @(11, "I", 224.5),
    @(11, "J", 1200),
    @(11, "K", 673.5),
    @(11, "L", 3600),
    @(11, "M", -533.5),
    @(11, "N", -3880),
    @(26, "H", 543.5),
    @(26, "I", 299.7143),
    @(26, "J", 2250),
    @(26, "K", 899.1428999999999),
    @(26, "L", 6750),
    @(26, "M", -611.1428999999999),
    @(26, "N", -7326),
    @(34, "H", 2807.5),
    @(34, "J", 2989),
    @(34, "L", 8967),
    @(34, "N", -9135),
    @(132, "H", 4176),
    @(132, "I", 2245.5),
    @(132, "K", 20209.5),
    @(132, "M", -17679.5),
    @(135, "H", 800.8),
    @(135, "I", 701.5),
    @(135, "J", 867),
    @(135, "K", 6313.5),
    @(135, "L", 7803),
    @(135, "M", -3778.5),
    @(135, "N", -12873)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @(126, "H", 2253.7144),
    @(126, "I", 1777),
    @(126, "K", 5331),
    @(126, "M", -2861),
    @(132, "H", 3334.9565),
    @(132, "I", 2774.2632),
    @(132, "J", 5998.25),
    @(132, "K", 8322.7896),
    @(132, "L", 17994.75),
    @(132, "M", -5792.7896),
    @(132, "N", -23054.75)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @(132, "H", 4612.1333),
    @(132, "I", 4283.5),
    @(132, "J", 4831.222),
    @(132, "K", 12850.5),
    @(132, "L", 14493.666),
    @(132, "M", -10320.5),
    @(132, "N", -19553.666),
    @(136, "H", 31758.934),
    @(136, "I", 10913.143),
    @(136, "K", 32739.429),
    @(136, "M", -30189.429)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @(132, "H", 3148.1428),
    @(132, "I", 2866.318),
    @(132, "J", 4181.5),
    @(132, "K", 8598.954000000002),
    @(132, "L", 12544.5),
    @(132, "M", -6068.954000000002),
    @(132, "N", -17604.5),
    @(136, "H", 8000.1),
    @(136, "I", 8000.1),
    @(136, "K", 24000.3),
    @(136, "M", -21450.3)
)
foreach ($u in $updates) {
    $row = $u[0]; $col = $u[1]; $val = $u[2]
    if ($null -eq $val) {
        $ws.Range("$col$row").ClearContents()
    } else {
        $ws.Range("$col$row").Value = $val
    }
}
